$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from LoginData to Sheet1
$ws.Name = "Sheet1"

# Header row (row 1) - replace old 7-column header with new 9-column header
$ws.Range("A1").Value = "testCase"
$ws.Range("B1").Value = "testType"
$ws.Range("C1").Value = "username"
$ws.Range("D1").Value = "password"
$ws.Range("E1").Value = "module"
$ws.Range("F1").Value = "expectedResult"
$ws.Range("G1").Value = "executeFlag"
$ws.Range("H1").Value = "environment"
$ws.Range("I1").Value = "priority`r"

# Data rows 2-18
# Row 2
$ws.Range("A2").Value = "TC501-1"
$ws.Range("B2").Value = "login"
$ws.Range("C2").Value = "Admin"
$ws.Range("D2").Value = "admin123"
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = "Login successful"
$ws.Range("G2").Value = $true
$ws.Range("H2").Value = "QA"
$ws.Range("I2").Value = "high`r"

# Row 3
$ws.Range("A3").Value = "TC501-2"
$ws.Range("B3").Value = "login"
$ws.Range("C3").Value = "testuser"
$ws.Range("D3").Value = "test123"
$ws.Range("E3").ClearContents()
$ws.Range("F3").Value = "Login successful"
$ws.Range("G3").Value = $true
$ws.Range("H3").Value = "QA"
$ws.Range("I3").Value = "high`r"

# Row 4
$ws.Range("A4").Value = "TC501-3"
$ws.Range("B4").Value = "login"
$ws.Range("C4").Value = "manager"
$ws.Range("D4").Value = "manager123"
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = "Login successful"
$ws.Range("G4").Value = $true
$ws.Range("H4").Value = "QA"
$ws.Range("I4").Value = "medium`r"

# Row 5
$ws.Range("A5").Value = "TC503-1"
$ws.Range("B5").Value = "navigation"
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("E5").Value = "Admin"
$ws.Range("F5").Value = "Admin page"
$ws.Range("G5").Value = $true
$ws.Range("H5").Value = "QA"
$ws.Range("I5").Value = "high`r"

# Row 6
$ws.Range("A6").Value = "TC503-2"
$ws.Range("B6").Value = "navigation"
$ws.Range("E6").Value = "PIM"
$ws.Range("F6").Value = "PIM page"
$ws.Range("G6").Value = $true
$ws.Range("H6").Value = "QA"
$ws.Range("I6").Value = "high`r"

# Row 7
$ws.Range("A7").Value = "TC503-3"
$ws.Range("B7").Value = "navigation"
$ws.Range("E7").Value = "Leave"
$ws.Range("F7").Value = "Leave page"
$ws.Range("G7").Value = $true
$ws.Range("H7").Value = "QA"
$ws.Range("I7").Value = "medium`r"

# Row 8
$ws.Range("A8").Value = "TC503-4"
$ws.Range("B8").Value = "navigation"
$ws.Range("E8").Value = "Time"
$ws.Range("F8").Value = "Time page"
$ws.Range("G8").Value = $true
$ws.Range("H8").Value = "QA"
$ws.Range("I8").Value = "medium`r"

# Row 9
$ws.Range("A9").Value = "TC503-5"
$ws.Range("B9").Value = "navigation"
$ws.Range("E9").Value = "Recruitment"
$ws.Range("F9").Value = "Recruitment page"
$ws.Range("G9").Value = $true
$ws.Range("H9").Value = "QA"
$ws.Range("I9").Value = "medium`r"

# Row 10
$ws.Range("A10").Value = "TC503-6"
$ws.Range("B10").Value = "navigation"
$ws.Range("E10").Value = "My Info"
$ws.Range("F10").Value = "My Info page"
$ws.Range("G10").Value = $true
$ws.Range("H10").Value = "QA"
$ws.Range("I10").Value = "low`r"

# Row 11
$ws.Range("A11").Value = "TC503-7"
$ws.Range("B11").Value = "navigation"
$ws.Range("E11").Value = "Performance"
$ws.Range("F11").Value = "Performance page"
$ws.Range("G11").Value = $true
$ws.Range("H11").Value = "QA"
$ws.Range("I11").Value = "low`r"

# Row 12
$ws.Range("A12").Value = "TC503-8"
$ws.Range("B12").Value = "navigation"
$ws.Range("E12").Value = "Dashboard"
$ws.Range("F12").Value = "Dashboard page"
$ws.Range("G12").Value = $true
$ws.Range("H12").Value = "QA"
$ws.Range("I12").Value = "high`r"

# Row 13
$ws.Range("A13").Value = "TC503-9"
$ws.Range("B13").Value = "navigation"
$ws.Range("E13").Value = "Directory"
$ws.Range("F13").Value = "Directory page"
$ws.Range("G13").Value = $true
$ws.Range("H13").Value = "QA"
$ws.Range("I13").Value = "low`r"

# Row 14
$ws.Range("A14").Value = "TC503-10"
$ws.Range("B14").Value = "navigation"
$ws.Range("E14").Value = "Maintenance"
$ws.Range("F14").Value = "Maintenance page"
$ws.Range("G14").Value = $true
$ws.Range("H14").Value = "QA"
$ws.Range("I14").Value = "low`r"

# Row 15
$ws.Range("A15").Value = "TC503-11"
$ws.Range("B15").Value = "navigation"
$ws.Range("E15").Value = "Buzz"
$ws.Range("F15").Value = "Buzz page"
$ws.Range("G15").Value = $true
$ws.Range("H15").Value = "QA"
$ws.Range("I15").Value = "low`r"

# Row 16
$ws.Range("A16").Value = "TC502-1"
$ws.Range("B16").Value = "menu-verify"
$ws.Range("E16").Value = "Admin"
$ws.Range("F16").Value = "Admin menu visible"
$ws.Range("G16").Value = $true
$ws.Range("H16").Value = "QA"
$ws.Range("I16").Value = "medium`r"

# Row 17
$ws.Range("A17").Value = "TC502-2"
$ws.Range("B17").Value = "menu-verify"
$ws.Range("E17").Value = "PIM"
$ws.Range("F17").Value = "PIM menu visible"
$ws.Range("G17").Value = $true
$ws.Range("H17").Value = "QA"
$ws.Range("I17").Value = "medium`r"

# Row 18
$ws.Range("A18").Value = "TC502-3"
$ws.Range("B18").Value = "menu-verify"
$ws.Range("E18").Value = "Leave"
$ws.Range("F18").Value = "Leave menu visible"
$ws.Range("G18").Value = $true
$ws.Range("H18").Value = "QA"
$ws.Range("I18").Value = "medium"
